$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 46005
$ws.Range("B5").Value = 48848

$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("B5").NumberFormat = $ws.Range("B4").NumberFormat
